# RESTORE: Recover all 973 original multi-industry template files from commit 168d9c4
#
# Relabels the "IT" themed training modules / roles back to the original
# AI/ML themed multi-industry template wording on the two schedule sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: "Training Schedule Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Schedule Overview")

$ws1.Range("A9").Value  = "AI/ML Fundamentals (AI-101)"
$ws1.Range("A10").Value = "AI/ML Platform Overview (AI-102)"
$ws1.Range("B11").Value = "Business Analysts"
$ws1.Range("B12").Value = "Data Scientists"
$ws1.Range("B13").Value = "ML Engineers, IT"
$ws1.Range("B14").Value = "ML Engineers, QA"

# ---------------------------------------------------------------------
# Sheet: "Detailed Training Schedule"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Detailed Training Schedule")

$ws2.Range("B4").Value = "AI/ML Fundamentals"
$ws2.Range("B5").Value = "AI/ML Platform Overview"

$ws2.Range("C6").Value  = "Business Analysts"
$ws2.Range("C7").Value  = "Business Analysts"
$ws2.Range("C8").Value  = "Business Analysts"
$ws2.Range("C9").Value  = "Data Scientists"
$ws2.Range("C10").Value = "Data Scientists"
$ws2.Range("C11").Value = "ML Engineers, IT"
$ws2.Range("C12").Value = "ML Engineers, IT"
$ws2.Range("C13").Value = "ML Engineers, QA"
$ws2.Range("C14").Value = "ML Engineers, QA"
